$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 32: AdaBoost -> GAN for Tabular Data
$ws.Range("D32").Value = "GAN for Tabular Data (Data Augmentation)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/334"

# Row 36: Active Learning in Semiconductor Manufacturing -> Intermediate Human Pose Estimation
$ws.Range("D36").Value = "Intermediate Human Pose Estimation"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/332"

# Row 41: AI game dev article -> Data Mesh metadata management article
$ws.Range("D41").Value = "Data Mesh 관점에서 효율적인 메타데이터 관리"
$ws.Range("E41").Value = "http://cloudinsight.net/data/data-mesh-%ea%b4%80%ec%a0%90%ec%97%90%ec%84%9c-%ed%9a%a8%ec%9c%a8%ec%a0%81%ec%9d%b8-%eb%a9%94%ed%83%80%eb%8d%b0%ec%9d%b4%ed%84%b0-%ea%b4%80%eb%a6%ac/"
